# "Generate Report for handoff"
#
# The localization-status report moves from "Handoff transform failed" to
# "Ready for handoff" on all three sheets, and the zh-cn / de-de detail
# sheets get a freshly generated handoff package: a link to the new .xlf
# target file, the handoff timestamp, and the dependency/handoff-reason
# flag flips from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$commitBase = "https://github.com/OpenLocalizationTest/oltest/blob/fc80edfd34bbf1a26328ea04480c93ae506e1594/e2e"
$baseName   = "3d029bba-3824-48aa-ba83-1438ac837909"
$revision   = "a7aa5fe93771e23b79815cc1d4b6756591e8bd8a"

# --- Status rolls from "Handoff transform failed" to "Ready for handoff" ---
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("B2").Value = "Ready for handoff"
$dede.Range("B2").Value = "Ready for handoff"

# --- zh-cn: newly produced handoff package ---
$zhcnFile = "$baseName.$revision.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "$commitBase/$zhcnFile", "", "", $zhcnFile)
$zhcn.Range("C2").Style = "HyperLink"
$zhcn.Range("C2").Font.Underline = 2
$zhcn.Range("C2").Font.Color = 15570276
$zhcn.Range("D2").Value = "2016-01-18 06:59:30"
$zhcn.Range("H2").Value = "Include"

# --- de-de: newly produced handoff package ---
$dedeFile = "$baseName.$revision.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("C2"), "$commitBase/$dedeFile", "", "", $dedeFile)
$dede.Range("C2").Style = "HyperLink"
$dede.Range("C2").Font.Underline = 2
$dede.Range("C2").Font.Color = 15570276
$dede.Range("D2").Value = "2016-01-18 06:59:39"
$dede.Range("H2").Value = "Include"
